# Update "想去人数" (want-to-go count) values in the F column
# for both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 565
$wsExhibit.Range("F10").Value = 6916
$wsExhibit.Range("F12").Value = 386
$wsExhibit.Range("F13").Value = 3287
$wsExhibit.Range("F14").Value = 225
$wsExhibit.Range("F15").Value = 404
$wsExhibit.Range("F17").Value = 570
$wsExhibit.Range("F18").Value = 40

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 565
$wsAll.Range("F13").Value = 6916
$wsAll.Range("F16").Value = 386
$wsAll.Range("F17").Value = 3287
$wsAll.Range("F18").Value = 225
$wsAll.Range("F19").Value = 404
$wsAll.Range("F21").Value = 570
$wsAll.Range("F22").Value = 40
